$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty template row (row 13) with a new
# "mejora" (improvement) request: Reporte de No conformidades.
$ws.Range("A13").Value = "Reporte de No conformidades"
$ws.Range("B13").Value = "En la plantilla de Reporte de no conformidades `n1. Quitar la columna de recomendaciones`n2. Agregar todos los procesos de la organización en la lista de la columna de proceso"
$ws.Range("C13").Value = "Aseguramiento de la calidad"
$ws.Range("D13").Value = "Actualización"
$ws.Range("E13").Value = "Aprobado "

# The description cell wraps, like the other description cells in the table.
$ws.Range("B13").WrapText = $true

# Tipo/Estatus columns pick up the same font styling already used by the
# rows above (D12:E12) instead of the font used by the blank template rows.
$ws.Range("D12:E12").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# Column A needs to be a bit wider to comfortably fit the new row's text.
$ws.Columns.Item(1).ColumnWidth = 33.1666666666666

# Leave the selection where the user ended up after adding the row.
$ws.Range("E12:E13").Select()
